$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 86.666664
$ws.Range("I5").Value = 86.666664
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 86.666664
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 28.333336
$ws.Range("N5").ClearContents()

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H64").Value = 1375
$ws.Range("I64").Value = 1000
$ws.Range("K64").Value = 1000
$ws.Range("M64").Value = -752

$ws.Range("H67").Value = 1375
$ws.Range("I67").Value = 1000
$ws.Range("K67").Value = 1000
$ws.Range("M67").Value = -142

$ws.Range("H76").Value = 4801.5
$ws.Range("I76").Value = 4801.5
$ws.Range("K76").Value = 4801.5
$ws.Range("M76").Value = -4486.5

$ws.Range("H79").Value = 4801.5
$ws.Range("I79").Value = 4801.5
$ws.Range("K79").Value = 4801.5
$ws.Range("M79").Value = -3709.5

$ws.Range("H92").Value = 602.5
$ws.Range("I92").Value = 771.6667
$ws.Range("J92").Value = 348.75
$ws.Range("K92").Value = 771.6667
$ws.Range("L92").Value = 348.75
$ws.Range("M92").Value = 476.3333
$ws.Range("N92").Value = -2844.75

$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492

$ws.Range("H100").Value = 966
$ws.Range("J100").Value = 1899
$ws.Range("L100").Value = 1899
$ws.Range("N100").Value = -2981

$ws.Range("H113").Value = 4000
$ws.Range("J113").Value = 4000
$ws.Range("L113").Value = 4000
$ws.Range("N113").Value = -10508

$ws.Range("H138").Value = 9743.333000000001
$ws.Range("I138").Value = 3397.3333
$ws.Range("J138").Value = 12916.333
$ws.Range("K138").Value = 10191.9999
$ws.Range("L138").Value = 38748.999
$ws.Range("M138").Value = -5051.999899999999
$ws.Range("N138").Value = -49028.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 3949.875
$ws.Range("I28").Value = 3949.875
$ws.Range("K28").Value = 3949.875
$ws.Range("M28").Value = -3757.875

$ws.Range("H32").Value = 759.8570999999999
$ws.Range("I32").Value = 759.8570999999999
$ws.Range("K32").Value = 759.8570999999999
$ws.Range("M32").Value = -472.8570999999999

$ws.Range("H99").Value = 3949.875
$ws.Range("I99").Value = 3949.875
$ws.Range("K99").Value = 3949.875
$ws.Range("M99").Value = -954.875

$ws.Range("H102").Value = 727.5
$ws.Range("I102").Value = 727.5
$ws.Range("K102").Value = 727.5
$ws.Range("M102").Value = 894.5

$ws.Range("H125").Value = 43725
$ws.Range("J125").Value = 43725
$ws.Range("L125").Value = 43725
$ws.Range("N125").Value = -53565

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 4497.5
$ws.Range("I75").Value = 4497.5
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 4497.5
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -3561.5
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 4497.5
$ws.Range("I78").Value = 4497.5
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 13492.5
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -8812.5
$ws.Range("N78").ClearContents()

$ws.Range("H86").Value = 1849.5
$ws.Range("I86").Value = 1849.5
$ws.Range("K86").Value = 1849.5
$ws.Range("M86").Value = -726.5

$ws.Range("H89").Value = 1849.5
$ws.Range("I89").Value = 1849.5
$ws.Range("K89").Value = 9247.5
$ws.Range("M89").Value = -3631.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 40.666668
$ws.Range("I7").Value = 32.666668
$ws.Range("K7").Value = 32.666668
$ws.Range("M7").Value = 80.333332

$ws.Range("H22").Value = 565
$ws.Range("I22").Value = 386.6
$ws.Range("J22").Value = 788
$ws.Range("K22").Value = 386.6
$ws.Range("L22").Value = 788
$ws.Range("M22").Value = -36.60000000000002
$ws.Range("N22").Value = -1488

$ws.Range("H132").Value = 3522.2144
$ws.Range("I132").Value = 1665.8182
$ws.Range("J132").Value = 10329
$ws.Range("K132").Value = 4997.4546
$ws.Range("L132").Value = 30987
$ws.Range("M132").Value = -2467.4546
$ws.Range("N132").Value = -36047

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 333567
$ws.Range("I16").Value = 333567
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1000701
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1000528
$ws.Range("N16").ClearContents()

$ws.Range("H39").Value = 700
$ws.Range("I39").Value = 700
$ws.Range("K39").Value = 2100
$ws.Range("M39").Value = -1806

$ws.Range("H61").Value = 400
$ws.Range("I61").Value = 400
$ws.Range("K61").Value = 1200
$ws.Range("M61").Value = -985

$ws.Range("H129").Value = 1504
$ws.Range("I129").Value = 515
$ws.Range("K129").Value = 1545
$ws.Range("M129").Value = 3455

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 14000
$ws.Range("J36").Value = 14000
$ws.Range("L36").Value = 14000
$ws.Range("N36").Value = -14970

$ws.Range("H40").Value = 10000
$ws.Range("J40").Value = 10000
$ws.Range("L40").Value = 10000
$ws.Range("N40").Value = -10302

$ws.Range("H41").Value = 1859.8
$ws.Range("I41").Value = 1074.75
$ws.Range("J41").Value = 5000
$ws.Range("K41").Value = 1074.75
$ws.Range("L41").Value = 5000
$ws.Range("M41").Value = -719.75
$ws.Range("N41").Value = -5710

$ws.Range("H44").Value = 33000
$ws.Range("I44").Value = 33000
$ws.Range("K44").Value = 33000
$ws.Range("M44").Value = -32404

$ws.Range("H46").Value = 6523
$ws.Range("J46").Value = 9046
$ws.Range("L46").Value = 9046
$ws.Range("N46").Value = -9358

$ws.Range("H92").Value = 1611.75
$ws.Range("J92").Value = 1611.75
$ws.Range("L92").Value = 1611.75
$ws.Range("N92").Value = -5355.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2398.5
$ws.Range("I100").Value = 2398
$ws.Range("J100").Value = 2399
$ws.Range("K100").Value = 2398
$ws.Range("L100").Value = 2399
$ws.Range("M100").Value = -1857
$ws.Range("N100").Value = -3481

$ws.Range("H104").Value = 8321.75
$ws.Range("J104").Value = 8321.75
$ws.Range("L104").Value = 8321.75
$ws.Range("N104").Value = -15309.75

